$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.205.35'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.59%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.861.67'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.66%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7142'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.05%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '240.26'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.31%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3085'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.11%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07685'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.85%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.14%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08324'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.88%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.960.83'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.12%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7164'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.02%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.212'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.06%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '90.70'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.52%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '29.328.78'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.30%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.946'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.89%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.172.85'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.07%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007800'
$ws.Range('D20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.90%  '

$ws.Range('E22').Value = '  +0.19%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.993'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.53%  '

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1611'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.81%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '162.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.900'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.03%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.44%  '

$ws.Range('E29').Value = '  -1.16%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.436'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.56%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.495'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.91%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.242'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.25%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05183'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.57%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7906'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +9.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.922'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.06%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.171'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.687'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.31%  '

$ws.Range('E38').Value = '  -0.15%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.693'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.72%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.183.37'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.30%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.244'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.91%  '

$ws.Range('E42').Value = '  -0.23%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.94'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.94%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.0000'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.068.56'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.48%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.39'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.84%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5205'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.56%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.777'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.24%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.340'
$ws.Range('D49').Style = 'Normal'

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.012'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.17%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.065'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.77%  '
